$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue 'D2' '97.008.18'
Set-TextValue 'E2' '  +0.10%  '
Set-TextValue 'D3' '3.714.71'
Set-TextValue 'E3' '  +0.47%  '
Set-TextValue 'E4' '  -0.02%  '
Set-TextValue 'D5' '237.10'
Set-TextValue 'E5' '  -3.29%  '
Set-TextValue 'E6' '  -0.81%  '
Set-TextValue 'D7' '652.70'
Set-TextValue 'E7' '  -2.56%  '
Set-TextValue 'E8' '  -0.29%  '
Set-TextValue 'E9' '  -0.02%  '
Set-TextValue 'E10' '  -5.55%  '
Set-TextValue 'D11' '3.712.86'
Set-TextValue 'E11' '  +0.54%  '
Set-TextValue 'B12' 'Avalanche'
Set-TextValue 'C12' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D12' '44.16'
Set-TextValue 'E12' '  -2.79%  '
Set-TextValue 'B13' 'ShibaInu'
Set-TextValue 'C13' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D13' '0.0000304'
Set-TextValue 'E13' '  +13.61%  '
Set-TextValue 'D14' '0.207'
Set-TextValue 'E14' '  +0.60%  '
Set-TextValue 'D15' '6.76'
Set-TextValue 'E15' '  +2.47%  '
Set-TextValue 'D16' '4.404.09'
Set-TextValue 'E16' '  +0.38%  '
Set-TextValue 'D17' '96.679.46'
Set-TextValue 'E17' '  -0.05%  '
Set-TextValue 'D18' '8.87'
Set-TextValue 'E18' '  -2.06%  '
Set-TextValue 'D19' '3.722.46'
Set-TextValue 'E19' '  +0.77%  '
Set-TextValue 'D20' '13.17'
Set-TextValue 'E20' '  +1.69%  '
Set-TextValue 'D21' '18.74'
Set-TextValue 'E21' '  +1.03%  '
Set-TextValue 'D22' '0.505'
Set-TextValue 'E22' '  -6.20%  '
Set-TextValue 'D23' '523.11'
Set-TextValue 'E23' '  +1.17%  '
Set-TextValue 'E24' '  -1.58%  '
Set-TextValue 'D25' '0.0000211'
Set-TextValue 'E25' '  +1.03%  '
Set-TextValue 'E26' '  -0.39%  '
Set-TextValue 'D27' '101.42'
Set-TextValue 'E27' '  -0.39%  '
Set-TextValue 'D28' '0.188'
Set-TextValue 'E28' '  +12.06%  '
Set-TextValue 'D29' '13.37'
Set-TextValue 'E29' '  +1.76%  '
Set-TextValue 'E30' '  -2.41%  '
Set-TextValue 'D31' '12.16'
Set-TextValue 'E31' '  -0.03%  '
Set-TextValue 'D32' '1.00'
Set-TextValue 'E32' '  +0.19%  '
Set-TextValue 'E33' '  +1.14%  '
Set-TextValue 'E34' '  +7.44%  '
Set-TextValue 'D35' '1.00'
Set-TextValue 'E35' '  +0.10%  '
Set-TextValue 'D36' '32.35'
Set-TextValue 'E36' '  -2.04%  '
Set-TextValue 'D37' '651.10'
Set-TextValue 'E37' '  +6.33%  '
Set-TextValue 'D38' '0.601'
Set-TextValue 'E38' '  +1.72%  '
Set-TextValue 'D39' '8.83'
Set-TextValue 'E39' '  -0.16%  '
Set-TextValue 'E40' '  +0.04%  '
Set-TextValue 'D41' '40.97'
Set-TextValue 'E41' '  -3.65%  '
Set-TextValue 'D42' '6.79'
Set-TextValue 'E42' '  +9.68%  '
Set-TextValue 'E43' '  +3.05%  '
Set-TextValue 'E44' '  -2.04%  '
Set-TextValue 'D45' '0.963'
Set-TextValue 'E45' '  -0.02%  '
Set-TextValue 'B46' 'Algorand'
Set-TextValue 'C46' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D46' '0.444'
Set-TextValue 'E46' '  +2.73%  '
Set-TextValue 'B47' 'VeChain'
Set-TextValue 'C47' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D47' '0.0454'
Set-TextValue 'E47' '  -0.73%  '
Set-TextValue 'D48' '2.29'
Set-TextValue 'E48' '  -0.82%  '
Set-TextValue 'E49' '  -0.01%  '
Set-TextValue 'E50' '  -1.26%  '
Set-TextValue 'E51' '  +0.42%  '
